$d = $word.ActiveDocument

# Locate the last existing "10 Checklist" bullet item so the new bullets are
# appended right after it (and before the final section break), inheriting
# the same ListParagraph style / numId 19 bullet numbering.
$rng = $d.Content
$found = $rng.Find.Execute("View Models for binding and responsive logic for what is done", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph for Change Log / Checklist insertion"
}

# Collapse to the end of that paragraph's text so InsertParagraphAfter()
# creates a new list paragraph right after it.
$rng.Collapse(0)

$newItems = @(
    "Got rid of the purple outline",
    "Added ability to add and see multiple interviews",
    "Changed the delete button design",
    "Added Navigation back to application page from details page",
    "Made some modifications to presentation on application details page based on user feedback"
)

foreach ($item in $newItems) {
    $rng.InsertParagraphAfter()
    $count = $d.Paragraphs.Count
    $newPara = $d.Paragraphs.Item($count)
    $newPara.Range.Text = $item
    $rng = $newPara.Range
    $rng.Collapse(0)
}

Write-Output "Inserted $($newItems.Count) checklist items."
